# Golang REST header compatibility
# Rename several REST header/body parameter names to Go-style hyphenated
# headers (e.g. UserName -> Username, ExpireSeconds -> Expire-Seconds, ...)
# and flatten the previously rich-text cells down to plain text, matching
# the authoritative commit diff. Also restores the sheet's frozen-pane /
# selection to the top of the data (B2) instead of the scrolled state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: POST /appmesh/login -----------------------------------------
$ws.Range("J2").Value = "curl -X POST --cacert /home/cent/ssl/out/centos/centos.bundle.crt -H ""Username:`$(echo -n admin | base64)"" -H ""password:`$(echo -n Admin123 | base64)"" https://localhost:6060/appmesh/login`ncurl -X POST -k -H ""Username:`$(echo -n admin | base64)"" -H ""password:`$(echo -n Admin123 | base64)"" https://localhost:6060/appmesh/login"
$ws.Range("E2").Value = "Required:`nUsername=base64(uname)`nPassword=base64(passwd)`nOptional:`nExpire-Seconds=600"

# --- Row 19: GET /appmesh/file/download ----------------------------------
$ws.Range("E19").Value = "Required：`nFile-Path=/opt/remote/filename"

# --- Row 20: POST /appmesh/file/upload -----------------------------------
$ws.Range("E20").Value = "Required：`nFile-Path=/opt/remote/filename`nOptional:`n  File-Mode=755`n  File-User=root"

# --- Row 3: POST /appmesh/auth -------------------------------------------
$ws.Range("E3").Value = "Required:`nAuthorization=""Bearer ""+access_token`nOptional:`nAuth-Permission=`${permission_key}"

# --- Restore frozen pane / selection to the top of the sheet -------------
# (was scrolled to show row 6 at top-left with G9 selected; move back to B2)
$ws.Range("B2").Select()
